# Journal de planification - ajout des entrées du 25 et 26.08.2023
# (Séance hebdomadaire avec le maître de diplôme + finalisation de la
# schématique Ethernet), et renommage de l'entrée du rendez-vous du 23.08.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Renomme l'ancien "Rendez-vous avec maitre de diplôme" (ligne 10) ---
$ws.Range("C10").Value = "Séance avec maitre de diplôme"

# --- 2. Nouvelle entrée du 25.08.2023 (ligne 15-18), fusion de la date ---
$ws.Range("A15").Value = 45163
$ws.Range("B15").Value = "08:00 - 10:00"
$ws.Range("C15").Value = "Séance hebdomadaire avec maître de diplôme"

$ws.Range("B16").Value = "10:00 - 12:00 "
$ws.Range("C16").Value = "Rédaction du procés-verbal de la séance hebdomadaire"

$ws.Range("B17").Value = "13:00 - 16:00"
$ws.Range("C17").Value = "Réalisation de la schématique, recherches de footprints et de composants."

$ws.Range("B18").Value = "16:00 - 17:00"
$ws.Range("C18").Value = "Réalisation de tests sur la consommation de courant des différents appareils. "

# --- 3. Nouvelle entrée du 26.08.2023 (ligne 19) ---
$ws.Range("A19").Value = 45164
$ws.Range("B19").Value = "10:00 - 16:00 "
$ws.Range("C19").Value = "Réalisation la schématique du connecteur et contrôleur Ethernet. Recherche du fonctionnement et des footprints des composants nécessaires"

# --- 4. Mise en forme de la colonne date (format date + centré), copiée
#        depuis une cellule de date existante pour rester cohérent avec le
#        reste du tableau (évite de dupliquer le format dans la feuille de
#        styles) ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A15:A18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null       # xlPasteFormats
$excel.CutCopyMode = $false

# --- 5. Fusion de la cellule de date pour le bloc du 25.08.2023 ---
$ws.Range("A15:A18").Merge() | Out-Null

# --- 6. Hauteur de la ligne 19 (texte sur deux lignes) ---
$ws.Rows.Item(19).RowHeight = 30

# --- 7. Sélection active ---
$ws.Range("D15").Select() | Out-Null
